{"js": "// Fix \"bug with empty notes\": a handful of footnotes had bad/placeholder\n// text (a stray trailing \"a\", or just a lone \"\u0f0d\" left over from a botched\n// split), and one trailing footnote reference at the very end of the\n// document was a leftover empty note that should never have been inserted.\n\nconst footnotes = context.document.body.footnotes;\nfootnotes.load(\"items\");\nawait context.sync();\n\n// Load every footnote's body text so we can find the ones that need fixing.\nfootnotes.items.forEach((fn) => fn.body.load(\"text\"));\nawait context.sync();\n\nconst TRAILING_A_MARK = \"\u0f54\u0f7a\u0f0b\u0f45\u0f72\u0f53\u0f0da\";\nconst isBlankNote = (text) => text.replace(/[\\u0002\\u0003\\s]/g, \"\") === \"\u0f0d\";\n\nlet noteWithStrayA = null;\nconst blankNotes = [];\n\nfor (const fn of footnotes.items) {\n  const text = fn.body.text;\n  if (text.indexOf(TRAILING_A_MARK) !== -1) {\n    noteWithStrayA = fn;\n  } else if (isBlankNote(text)) {\n    blankNotes.push(fn);\n  }\n}\n\n// The very last footnote in the document is the stray, empty note that was\n// accidentally left behind; it must be removed entirely (reference + note).\nconst lastFootnote = footnotes.items[footnotes.items.length - 1];\n\n// 1) Strip the stray trailing \"a\" that leaked into this footnote's text.\nif (noteWithStrayA) {\n  noteWithStrayA.body.paragraphs.load(\"items\");\n  await context.sync();\n  const paragraph = noteWithStrayA.body.paragraphs.items[0];\n  const fixedText = paragraph.text\n    .replace(/\\u0002/g, \"\")\n    .replace(/\\r$/, \"\")\n    .replace(/a$/, \"\");\n  paragraph.insertText(fixedText, Word.InsertLocation.replace);\n}\n\n// 2) The first two blank notes (just \"\u0f0d\") get their real content restored.\n//    They are filled in document order, skipping the trailing one that will\n//    be deleted outright.\nconst replacementTexts = [\n  \" \u0f40\u0fb1\u0f72\u0f0d \u0f5e\u0f7a\u0f66\u0f0b\u0f54\u0f62\u0f0b\u0f58\u0f0b\u0f42\u0f5e\u0f53\u0f0b\u0f53\u0f44\u0f0b\u0f58\u0f7a\u0f51\u0f0d\",\n  \" \u0f58\u0f0d \u0f5e\u0f7a\u0f66\u0f0b\u0f54\u0f62\u0f0b\u0f58\u0f0b\u0f42\u0f5e\u0f53\u0f0b\u0f53\u0f44\u0f0b\u0f58\u0f7a\u0f51\u0f0d\",\n];\n\nconst blankNotesToFix = blankNotes.filter((fn) => fn !== lastFootnote);\nfor (let i = 0; i < blankNotesToFix.length && i < replacementTexts.length; i++) {\n  const fn = blankNotesToFix[i];\n  fn.body.paragraphs.load(\"items\");\n  await context.sync();\n  const paragraph = fn.body.paragraphs.items[0];\n  paragraph.insertText(replacementTexts[i], Word.InsertLocation.replace);\n}\n\nawait context.sync();\n\n// 3) Remove the trailing empty footnote completely (its in-text reference\n//    mark as well as its note definition).\nlastFootnote.reference.delete();\n\nawait context.sync();\n", "ps1": "# Fix \"bug with empty notes\": a handful of footnotes had bad/placeholder\n# text (a stray trailing \"a\", or just a lone \"\u0f0d\" left over from a botched\n# split), and one trailing footnote reference at the very end of the\n# document was a leftover empty note that should never have been inserted.\n\n$d = $word.ActiveDocument\n\n$strayAIndex = $null\n$blankIndexes = @()\n\nfor ($i = 1; $i -le $d.Footnotes.Count; $i++) {\n    $fn = $d.Footnotes.Item($i)\n    $txt = $fn.Range.Text\n    $clean = $txt -replace \"[\\u0002\\u0003\\s]\", \"\"\n    if ($txt -match \"a$\") {\n        $strayAIndex = $i\n    } elseif ($clean -eq [char]0x0F0D) {\n        $blankIndexes += $i\n    }\n}\n\n$lastIndex = $d.Footnotes.Count\n\n# 1) Strip the stray trailing \"a\" that leaked into this footnote's text.\nif ($strayAIndex) {\n    $fn = $d.Footnotes.Item($strayAIndex)\n    $fn.Range.Text = ($fn.Range.Text -replace \"a$\", \"\")\n}\n\n# 2) The first two blank notes (just \"\u0f0d\") get their real content restored.\n#    They are filled in document order, skipping the trailing one that will\n#    be deleted outright.\n$fixIndexes = $blankIndexes | Where-Object { $_ -ne $lastIndex }\n$replacementTexts = @(\" \u0f40\u0fb1\u0f72\u0f0d \u0f5e\u0f7a\u0f66\u0f0b\u0f54\u0f62\u0f0b\u0f58\u0f0b\u0f42\u0f5e\u0f53\u0f0b\u0f53\u0f44\u0f0b\u0f58\u0f7a\u0f51\u0f0d\", \" \u0f58\u0f0d \u0f5e\u0f7a\u0f66\u0f0b\u0f54\u0f62\u0f0b\u0f58\u0f0b\u0f42\u0f5e\u0f53\u0f0b\u0f53\u0f44\u0f0b\u0f58\u0f7a\u0f51\u0f0d\")\nfor ($j = 0; $j -lt $fixIndexes.Count -and $j -lt $replacementTexts.Count; $j++) {\n    $fn = $d.Footnotes.Item($fixIndexes[$j])\n    $fn.Range.Text = $replacementTexts[$j]\n}\n\n# 3) Remove the trailing empty footnote completely (its in-text reference\n#    mark as well as its note definition).\n$lastFootnote = $d.Footnotes.Item($lastIndex)\n$lastFootnote.Delete()\n"}
